# arreglo de mapas segun nivel
#
# Locates the two list paragraphs about the player avatar / username HUD
# and rewrites them so that:
#   - the "Se debe capturar un avatar..." paragraph's highlight becomes
#     green (was yellow) and its trailing " (falta imagen)" run is
#     reduced to just a single space;
#   - the "Cuando el usuario ingrese..." paragraph loses its trailing
#     " (falta imagen)" run entirely;
#   - the _GoBack bookmark moves from the end of the first paragraph to
#     the end of the second paragraph.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Find-ParagraphContaining($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

$pAvatar = Find-ParagraphContaining $d "Se debe capturar un avatar"
if ($null -eq $pAvatar) {
    throw "Could not find the 'Se debe capturar un avatar' paragraph"
}

$pUsername = Find-ParagraphContaining $d "Cuando el usuario ingrese"
if ($null -eq $pUsername) {
    throw "Could not find the 'Cuando el usuario ingrese' paragraph"
}

$rAvatar = $pAvatar.Range
$rUsername = $pUsername.Range

$avatarXml = '<w:p ' + $wNs + '>' + `
    '<w:pPr>' + `
        '<w:pStyle w:val="Prrafodelista"/>' + `
        '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' + `
        '<w:rPr><w:highlight w:val="green"/></w:rPr>' + `
    '</w:pPr>' + `
    '<w:r>' + `
        '<w:rPr><w:highlight w:val="green"/></w:rPr>' + `
        '<w:t>Se debe capturar un avatar para los jugadores.</w:t>' + `
    '</w:r>' + `
    '<w:r>' + `
        '<w:rPr><w:highlight w:val="green"/></w:rPr>' + `
        '<w:t xml:space="preserve"> </w:t>' + `
    '</w:r>' + `
'</w:p>'

$usernameXml = '<w:p ' + $wNs + '>' + `
    '<w:pPr>' + `
        '<w:pStyle w:val="Prrafodelista"/>' + `
        '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' + `
    '</w:pPr>' + `
    '<w:r>' + `
        '<w:rPr><w:highlight w:val="green"/></w:rPr>' + `
        '<w:t xml:space="preserve">Cuando el usuario ingrese debe verse el </w:t>' + `
    '</w:r>' + `
    '<w:r>' + `
        '<w:rPr><w:highlight w:val="green"/></w:rPr>' + `
        '<w:t>username</w:t>' + `
    '</w:r>' + `
    '<w:r>' + `
        '<w:rPr><w:highlight w:val="green"/></w:rPr>' + `
        '<w:t xml:space="preserve"> del jugador</w:t>' + `
    '</w:r>' + `
    '<w:r>' + `
        '<w:rPr><w:highlight w:val="green"/></w:rPr>' + `
        '<w:t>, el nivel (en una imagen)</w:t>' + `
    '</w:r>' + `
    '<w:r>' + `
        '<w:rPr><w:highlight w:val="green"/></w:rPr>' + `
        '<w:t>,</w:t>' + `
    '</w:r>' + `
    '<w:r>' + `
        '<w:t xml:space="preserve"> </w:t>' + `
    '</w:r>' + `
    '<w:r>' + `
        '<w:rPr><w:highlight w:val="green"/></w:rPr>' + `
        '<w:t>cantidad de puntos</w:t>' + `
    '</w:r>' + `
    '<w:r>' + `
        '<w:t xml:space="preserve"> </w:t>' + `
    '</w:r>' + `
    '<w:r>' + `
        '<w:rPr><w:highlight w:val="yellow"/></w:rPr>' + `
        '<w:t>y el mundo en el que se encuentra</w:t>' + `
    '</w:r>' + `
    '<w:r>' + `
        '<w:t>.</w:t>' + `
    '</w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
'</w:p>'

# Replace the second paragraph first, then the first; either order is
# fine because Range objects captured before the edits keep tracking
# their own paragraph even as the document shifts.
$rUsername.InsertXML($usernameXml)
$rAvatar.InsertXML($avatarXml)
